$wb = $excel.ActiveWorkbook

# --- Pool sheet: add VL / SL columns (AP:AQ) ---
$pool = $wb.Worksheets.Item("Pool")
$pool.Activate()

$pool.Range("AP1").Value = "VL"
$pool.Range("AQ1").Value = "SL"
$pool.Range("AP2").Value = "15"
$pool.Range("AQ2").Value = "15"

# Copy the newly added range, mirroring the author's workflow (leaves a
# lingering marching-ants selection referenced from other sheets' views).
$pool.Range("AP1:AQ2").Copy()
$pool.Range("AP1:AQ2").Select()

# --- Touch the other sheets so their saved selection reflects the copy ---
$req = $wb.Worksheets.Item("Requirement")
$req.Activate()
$req.Range("F2").Select()

$sched = $wb.Worksheets.Item("Schedule")
$sched.Activate()
$sched.Range("A3").Select()

$hol = $wb.Worksheets.Item("hol")
$hol.Activate()
$hol.Range("I19").Select()

$wage = $wb.Worksheets.Item("Wage")
$wage.Activate()
$wage.Range("C3").Select()

# Return to Pool as the active sheet (matches tabSelected state in the target file)
$pool.Activate()
